$wb = $excel.ActiveWorkbook

# Sheet 1: 展览
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 244
$ws1.Range("F5").Value = 5
$ws1.Range("F6").Value = 2028
$ws1.Range("F8").Value = 650
$ws1.Range("F9").Value = 21
$ws1.Range("F10").Value = 161
$ws1.Range("G10").Value = 78
$ws1.Range("F11").Value = 144
$ws1.Range("F12").Value = 638
$ws1.Range("F13").Value = 31
$ws1.Range("F14").Value = 79
$ws1.Range("F15").Value = 1083
$ws1.Range("F18").Value = 180

# Sheet 2: 演出
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F6").Value = 9
$ws2.Range("F11").Value = 27
$ws2.Range("F13").Value = 7
$ws2.Range("F20").Value = 45

# Sheet 3: 本地生活
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 6315
$ws3.Range("F3").Value = 783
$ws3.Range("F4").Value = 1967
$ws3.Range("F5").Value = 193

# Sheet 4: 全部类型
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 6315
$ws4.Range("F3").Value = 783
$ws4.Range("F4").Value = 1967
$ws4.Range("F6").Value = 193
$ws4.Range("F12").Value = 244
$ws4.Range("F13").Value = 5
$ws4.Range("F14").Value = 9
$ws4.Range("F17").Value = 2028
$ws4.Range("F21").Value = 27
$ws4.Range("F22").Value = 650
$ws4.Range("F23").Value = 21
$ws4.Range("F24").Value = 161
$ws4.Range("G24").Value = 78
$ws4.Range("F26").Value = 144
$ws4.Range("F27").Value = 638
$ws4.Range("F28").Value = 31
$ws4.Range("F29").Value = 79
$ws4.Range("F30").Value = 7
$ws4.Range("F31").Value = 1083
$ws4.Range("F36").Value = 180
$ws4.Range("F41").Value = 45
